# Auto-generated edit script applying crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($sheet, $addr, $val) {
    $r = $sheet.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue $ws "D2" "29.851.92"
$ws.Range("E2").Value = "  -0.13%  "

Set-TextValue $ws "D3" "1.639.08"
$ws.Range("E3").Value = "  +0.93%  "

$ws.Range("E4").Value = "  +0.48%  "

Set-TextValue $ws "D5" "215.34"
$ws.Range("E5").Value = "  +0.55%  "

$ws.Range("E6").Value = "  -0.28%  "

$ws.Range("E7").Value = "  +0.51%  "

Set-TextValue $ws "D8" "28.91"
$ws.Range("E8").Value = "  -2.66%  "

$ws.Range("E9").Value = "  +0.74%  "

Set-TextValue $ws "D10" "0.0609"
$ws.Range("E10").Value = "  -0.18%  "

Set-TextValue $ws "D11" "0.0898"
$ws.Range("E11").Value = "  -1.54%  "

Set-TextValue $ws "D12" "1.873.25"
$ws.Range("E12").Value = "  +0.84%  "

Set-TextValue $ws "D13" "1.642.15"
$ws.Range("E13").Value = "  +1.34%  "

Set-TextValue $ws "D14" "0.590"
$ws.Range("E14").Value = "  +3.82%  "

Set-TextValue $ws "D15" "9.45"
$ws.Range("E15").Value = "  +7.45%  "

$ws.Range("E16").Value = "  -0.46%  "

Set-TextValue $ws "D17" "29.849.29"
$ws.Range("E17").Value = "  -0.27%  "

Set-TextValue $ws "D18" "64.42"
$ws.Range("E18").Value = "  -0.37%  "

Set-TextValue $ws "D19" "238.98"
$ws.Range("E19").Value = "  -2.16%  "

Set-TextValue $ws "D20" "0.0₃0703"
$ws.Range("E20").Value = "  -0.12%  "

$ws.Range("E21").Value = "  +0.47%  "

Set-TextValue $ws "D22" "9.95"
$ws.Range("E22").Value = "  +3.37%  "

$ws.Range("E23").Value = "  +0.34%  "

$ws.Range("E24").Value = "  +2.26%  "

Set-TextValue $ws "D25" "156.87"
$ws.Range("E25").Value = "  -0.12%  "

Set-TextValue $ws "D26" "15.55"
$ws.Range("E26").Value = "  -0.80%  "

$ws.Range("E27").Value = "  -1.16%  "

Set-TextValue $ws "D28" "6.62"
$ws.Range("E28").Value = "  +0.29%  "

$ws.Range("E29").Value = "  +0.42%  "

$ws.Range("E30").Value = "  +1.18%  "

$ws.Range("E31").Value = "  -0.65%  "

Set-TextValue $ws "D32" "3.38"
$ws.Range("E32").Value = "  +1.07%  "

$ws.Range("E33").Value = "  -0.96%  "

Set-TextValue $ws "D34" "1.420.13"

$ws.Range("E35").Value = "  +2.46%  "

Set-TextValue $ws "D36" "1.02"
$ws.Range("E36").Value = "  -1.33%  "

Set-TextValue $ws "D37" "2.70"
$ws.Range("E37").Value = "  -5.83%  "

$ws.Range("E38").Value = "  +1.65%  "

$ws.Range("E39").Value = "  +0.38%  "

Set-TextValue $ws "D40" "76.42"
$ws.Range("E40").Value = "  +10.32%  "

Set-TextValue $ws "D41" "0.564"
$ws.Range("E41").Value = "  +1.21%  "

$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws "D42" "0.833"
$ws.Range("E42").Value = "  -0.19%  "

$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws "D43" "0.0500"
$ws.Range("E43").Value = "  -1.59%  "

$ws.Range("E44").Value = "  -3.02%  "

$ws.Range("E45").Value = "  +0.56%  "

Set-TextValue $ws "D46" "1.00"
$ws.Range("E46").Value = "  -2.17%  "

Set-TextValue $ws "D47" "1.781.54"
$ws.Range("E47").Value = "  +0.88%  "

$ws.Range("B48").Value = "BitcoinSV"
$ws.Range("C48").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
Set-TextValue $ws "D48" "49.61"
$ws.Range("E48").Value = "  -8.19%  "

$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws "D49" "5.33"
$ws.Range("E49").Value = "  -1.75%  "

Set-TextValue $ws "D50" "93.17"
$ws.Range("E50").Value = "  +5.39%  "

$ws.Range("E51").Value = "  +1.45%  "
